$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the selection that Excel had saved in the sheetView so the sheet
# reopens without an active selection/highlight (matches re-run export).
$ws.Range("B1").Select()

# Row 3
$ws.Range("C3").Value = 2268578
$ws.Range("D3").Value = 2261252
$ws.Range("E3").Value = 2268644
$ws.Range("F3").Value = 2261318

# Row 4
$ws.Range("C4").Value = 4823533
$ws.Range("D4").Value = 4809213
$ws.Range("E4").Value = 6400938
$ws.Range("F4").Value = 6388496

# Row 5
$ws.Range("C5").Value = 2178867
$ws.Range("D5").Value = 2178867

# Row 6
$ws.Range("C6").Value = 2616
$ws.Range("D6").Value = 2469
$ws.Range("E6").Value = 2472
$ws.Range("F6").Value = 2315

# Row 7
$ws.Range("C7").Value = 2771
$ws.Range("D7").Value = 2507
$ws.Range("E7").Value = 64953
$ws.Range("F7").Value = 48472

# Row 9
$ws.Range("D9").Value = 20030

# Row 10
$ws.Range("C10").Value = 2268578
$ws.Range("D10").Value = 2261252
$ws.Range("E10").Value = 2268644
$ws.Range("F10").Value = 2261318

# Row 11
$ws.Range("D11").Value = 20030

# Row 12
$ws.Range("C12").Value = 20030
$ws.Range("D12").Value = 20030

# Row 13
$ws.Range("D13").Value = 40060

# Row 14
$ws.Range("E14").Value = 12
$ws.Range("F14").Value = 12

# Row 15
$ws.Range("E15").Value = 12
$ws.Range("F15").Value = 12

# Row 16
$ws.Range("C16").Value = 555220992
$ws.Range("D16").Value = 555745280
$ws.Range("E16").Value = 2337800192
$ws.Range("F16").Value = 2285895680
